# Auto-generated Excel COM-interop script
# Applies per-cell value updates to columns H:N across multiple sheets
# (market-price / profit recalculation snapshot), matching the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 509.26666  # H28: 542.6429000000001 -> 509.26666
$ws.Cells.Item(28, 9).Value = 402.7857  # I28: 430.53845 -> 402.7857
$ws.Cells.Item(28, 11).Value = 402.7857  # K28: 430.53845 -> 402.7857
$ws.Cells.Item(28, 13).Value = 82.21429999999998  # M28: 54.46154999999999 -> 82.21429999999998
$ws.Cells.Item(98, 8).Value = 3860.077  # H98: 3379.1333 -> 3860.077
$ws.Cells.Item(98, 9).Value = 3672.2856  # I98: 3275.75 -> 3672.2856
$ws.Cells.Item(98, 10).Value = 4079.1667  # J98: 3497.2856 -> 4079.1667
$ws.Cells.Item(98, 11).Value = 3672.2856  # K98: 3275.75 -> 3672.2856
$ws.Cells.Item(98, 12).Value = 4079.1667  # L98: 3497.2856 -> 4079.1667
$ws.Cells.Item(98, 13).Value = -2174.2856  # M98: -1777.75 -> -2174.2856
$ws.Cells.Item(98, 14).Value = -7075.1667  # N98: -6493.2856 -> -7075.1667
$ws.Cells.Item(115, 8).Value = 264.5  # H115: 243 -> 264.5
$ws.Cells.Item(115, 9).Value = 264.5  # I115: 243 -> 264.5
$ws.Cells.Item(115, 11).Value = 793.5  # K115: 729 -> 793.5
$ws.Cells.Item(115, 13).Value = 773.5  # M115: 838 -> 773.5
$ws.Cells.Item(122, 8).Value = 3860.077  # H122: 3379.1333 -> 3860.077
$ws.Cells.Item(122, 9).Value = 3672.2856  # I122: 3275.75 -> 3672.2856
$ws.Cells.Item(122, 10).Value = 4079.1667  # J122: 3497.2856 -> 4079.1667
$ws.Cells.Item(122, 11).Value = 11016.8568  # K122: 9827.25 -> 11016.8568
$ws.Cells.Item(122, 12).Value = 12237.5001  # L122: 10491.8568 -> 12237.5001
$ws.Cells.Item(122, 13).Value = -8566.856800000001  # M122: -7377.25 -> -8566.856800000001
$ws.Cells.Item(122, 14).Value = -17137.5001  # N122: -15391.8568 -> -17137.5001
$ws.Cells.Item(127, 8).Value = 5666.6665  # H127: 3718.7144 -> 5666.6665
$ws.Cells.Item(127, 9).Value = 4500  # I127: 2806.2 -> 4500
$ws.Cells.Item(127, 10).Value = 6250  # J127: 6000 -> 6250
$ws.Cells.Item(127, 11).Value = 13500  # K127: 8418.599999999999 -> 13500
$ws.Cells.Item(127, 12).Value = 18750  # L127: 18000 -> 18750
$ws.Cells.Item(127, 13).Value = -8540  # M127: -3458.599999999999 -> -8540
$ws.Cells.Item(127, 14).Value = -28670  # N127: -27920 -> -28670
$ws.Cells.Item(132, 8).Value = 2912.7727  # H132: 2738.1277 -> 2912.7727
$ws.Cells.Item(132, 9).Value = 2880.8647  # I132: 2678.05 -> 2880.8647
$ws.Cells.Item(132, 11).Value = 8642.5941  # K132: 8034.150000000001 -> 8642.5941
$ws.Cells.Item(132, 13).Value = -6112.5941  # M132: -5504.150000000001 -> -6112.5941
$ws.Cells.Item(135, 8).Value = 1342  # H135: 1198.2222 -> 1342
$ws.Cells.Item(135, 9).Value = 347.8  # I135: 297.83334 -> 347.8
$ws.Cells.Item(135, 11).Value = 3130.2  # K135: 2680.50006 -> 3130.2
$ws.Cells.Item(135, 13).Value = -595.2000000000003  # M135: -145.5000600000003 -> -595.2000000000003
$ws.Cells.Item(137, 8).Value = 2092.0588  # H137: 2081.611 -> 2092.0588
$ws.Cells.Item(137, 9).Value = 1529.2307  # I137: 1556 -> 1529.2307
$ws.Cells.Item(137, 11).Value = 4587.6921  # K137: 4668 -> 4587.6921
$ws.Cells.Item(137, 13).Value = -2037.6921  # M137: -2118 -> -2037.6921
$ws.Cells.Item(138, 8).Value = 2553.1  # H138: 2644.4092 -> 2553.1
$ws.Cells.Item(138, 9).Value = 1170.25  # I138: 1256.875 -> 1170.25
$ws.Cells.Item(138, 10).Value = 3475  # J138: 3437.2856 -> 3475
$ws.Cells.Item(138, 11).Value = 3510.75  # K138: 3770.625 -> 3510.75
$ws.Cells.Item(138, 12).Value = 10425  # L138: 10311.8568 -> 10425
$ws.Cells.Item(138, 13).Value = 1629.25  # M138: 1369.375 -> 1629.25
$ws.Cells.Item(138, 14).Value = -20705  # N138: -20591.8568 -> -20705

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8372.143  # H32: 8536.617 -> 8372.143
$ws.Cells.Item(32, 9).Value = 7147.794  # I32: 7280.1514 -> 7147.794
$ws.Cells.Item(32, 11).Value = 7147.794  # K32: 7280.1514 -> 7147.794
$ws.Cells.Item(32, 13).Value = -6860.794  # M32: -6993.1514 -> -6860.794
$ws.Cells.Item(45, 8).Value = 2914.182  # H45: 2567.125 -> 2914.182
$ws.Cells.Item(45, 9).Value = 2542.6667  # I45: 2433.8572 -> 2542.6667
$ws.Cells.Item(45, 10).Value = 3360  # J45: 3500 -> 3360
$ws.Cells.Item(45, 11).Value = 2542.6667  # K45: 2433.8572 -> 2542.6667
$ws.Cells.Item(45, 12).Value = 3360  # L45: 3500 -> 3360
$ws.Cells.Item(45, 13).Value = -2165.6667  # M45: -2056.8572 -> -2165.6667
$ws.Cells.Item(45, 14).Value = -4114  # N45: -4254 -> -4114
$ws.Cells.Item(61, 8).Value = 1680.4445  # H61: 1548.909 -> 1680.4445
$ws.Cells.Item(61, 9).Value = 1680.4445  # I61: 1602.4 -> 1680.4445
$ws.Cells.Item(61, 10).Value = 0  # J61: 1014 -> 0
$ws.Cells.Item(61, 11).Value = 1680.4445  # K61: 1602.4 -> 1680.4445
$ws.Cells.Item(61, 12).Value = 0  # L61: 1014 -> 0
$ws.Cells.Item(61, 13).Value = -1468.4445  # M61: -1390.4 -> -1468.4445
$ws.Cells.Item(61, 14).ClearContents()  # N61: -1438 -> (cleared)
$ws.Cells.Item(122, 8).Value = 2864.9092  # H122: 1920.8182 -> 2864.9092
$ws.Cells.Item(122, 9).Value = 2951.4  # I122: 1917.0476 -> 2951.4
$ws.Cells.Item(122, 11).Value = 8854.200000000001  # K122: 5751.142800000001 -> 8854.200000000001
$ws.Cells.Item(122, 13).Value = -6404.200000000001  # M122: -3301.142800000001 -> -6404.200000000001
$ws.Cells.Item(132, 8).Value = 3417.3809  # H132: 3540.75 -> 3417.3809
$ws.Cells.Item(132, 9).Value = 3038.25  # I132: 3148.158 -> 3038.25
$ws.Cells.Item(132, 11).Value = 9114.75  # K132: 9444.474 -> 9114.75
$ws.Cells.Item(132, 13).Value = -6584.75  # M132: -6914.474 -> -6584.75
$ws.Cells.Item(136, 8).Value = 1680.4445  # H136: 1548.909 -> 1680.4445
$ws.Cells.Item(136, 9).Value = 1680.4445  # I136: 1602.4 -> 1680.4445
$ws.Cells.Item(136, 10).Value = 0  # J136: 1014 -> 0
$ws.Cells.Item(136, 11).Value = 5041.333500000001  # K136: 4807.200000000001 -> 5041.333500000001
$ws.Cells.Item(136, 12).Value = 0  # L136: 3042 -> 0
$ws.Cells.Item(136, 13).Value = -2491.333500000001  # M136: -2257.200000000001 -> -2491.333500000001
$ws.Cells.Item(136, 14).ClearContents()  # N136: -8142 -> (cleared)

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2319.75  # H20: 2283.923 -> 2319.75
$ws.Cells.Item(20, 9).Value = 2026.7142  # I20: 1970.7778 -> 2026.7142
$ws.Cells.Item(20, 10).Value = 2730  # J20: 2988.5 -> 2730
$ws.Cells.Item(20, 11).Value = 2026.7142  # K20: 1970.7778 -> 2026.7142
$ws.Cells.Item(20, 12).Value = 2730  # L20: 2988.5 -> 2730
$ws.Cells.Item(20, 13).Value = -1779.7142  # M20: -1723.7778 -> -1779.7142
$ws.Cells.Item(20, 14).Value = -3224  # N20: -3482.5 -> -3224
$ws.Cells.Item(134, 8).Value = 8855.333000000001  # H134: 11139.6 -> 8855.333000000001
$ws.Cells.Item(134, 9).Value = 8855.333000000001  # I134: 12324.5 -> 8855.333000000001
$ws.Cells.Item(134, 10).Value = 0  # J134: 6400 -> 0
$ws.Cells.Item(134, 11).Value = 26565.999  # K134: 36973.5 -> 26565.999
$ws.Cells.Item(134, 12).Value = 0  # L134: 19200 -> 0
$ws.Cells.Item(134, 13).Value = -24030.999  # M134: -34438.5 -> -24030.999
$ws.Cells.Item(134, 14).ClearContents()  # N134: -24270 -> (cleared)

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 319.83334  # H19: 403.25 -> 319.83334
$ws.Cells.Item(19, 9).Value = 129  # I19: 146.71428 -> 129
$ws.Cells.Item(19, 10).Value = 701.5  # J19: 762.4 -> 701.5
$ws.Cells.Item(19, 11).Value = 129  # K19: 146.71428 -> 129
$ws.Cells.Item(19, 12).Value = 701.5  # L19: 762.4 -> 701.5
$ws.Cells.Item(19, 13).Value = 41  # M19: 23.28572 -> 41
$ws.Cells.Item(19, 14).Value = -1041.5  # N19: -1102.4 -> -1041.5
$ws.Cells.Item(24, 8).Value = 319.83334  # H24: 403.25 -> 319.83334
$ws.Cells.Item(24, 9).Value = 129  # I24: 146.71428 -> 129
$ws.Cells.Item(24, 10).Value = 701.5  # J24: 762.4 -> 701.5
$ws.Cells.Item(24, 11).Value = 129  # K24: 146.71428 -> 129
$ws.Cells.Item(24, 12).Value = 701.5  # L24: 762.4 -> 701.5
$ws.Cells.Item(24, 13).Value = 41  # M24: 23.28572 -> 41
$ws.Cells.Item(24, 14).Value = -1041.5  # N24: -1102.4 -> -1041.5
$ws.Cells.Item(42, 8).Value = 0  # H42: 500 -> 0
$ws.Cells.Item(42, 9).Value = 0  # I42: 500 -> 0
$ws.Cells.Item(42, 11).Value = 0  # K42: 500 -> 0
$ws.Cells.Item(42, 13).ClearContents()  # M42: 93 -> (cleared)
$ws.Cells.Item(58, 9).Value = 6449.8335  # I58: 5742.7144 -> 6449.8335
$ws.Cells.Item(58, 10).Value = 1599.875  # J58: 1614.1428 -> 1599.875
$ws.Cells.Item(58, 11).Value = 6449.8335  # K58: 5742.7144 -> 6449.8335
$ws.Cells.Item(58, 12).Value = 1599.875  # L58: 1614.1428 -> 1599.875
$ws.Cells.Item(58, 13).Value = -6246.8335  # M58: -5539.7144 -> -6246.8335
$ws.Cells.Item(58, 14).Value = -2005.875  # N58: -2020.1428 -> -2005.875
$ws.Cells.Item(87, 8).Value = 0  # H87: 10000 -> 0
$ws.Cells.Item(87, 10).Value = 0  # J87: 10000 -> 0
$ws.Cells.Item(87, 12).Value = 0  # L87: 10000 -> 0
$ws.Cells.Item(87, 14).ClearContents()  # N87: -12372 -> (cleared)
$ws.Cells.Item(90, 8).Value = 0  # H90: 10000 -> 0
$ws.Cells.Item(90, 10).Value = 0  # J90: 10000 -> 0
$ws.Cells.Item(90, 12).Value = 0  # L90: 30000 -> 0
$ws.Cells.Item(90, 14).ClearContents()  # N90: -41856 -> (cleared)
$ws.Cells.Item(132, 8).Value = 3199  # H132: 1745.7273 -> 3199
$ws.Cells.Item(132, 9).Value = 3199  # I132: 1745.7273 -> 3199
$ws.Cells.Item(132, 11).Value = 9597  # K132: 5237.1819 -> 9597
$ws.Cells.Item(132, 13).Value = -7067  # M132: -2707.1819 -> -7067
$ws.Cells.Item(134, 8).Value = 1525.125  # H134: 1532.375 -> 1525.125
$ws.Cells.Item(134, 9).Value = 1652.579  # I134: 1661.7368 -> 1652.579
$ws.Cells.Item(134, 11).Value = 4957.737  # K134: 4985.2104 -> 4957.737
$ws.Cells.Item(134, 13).Value = -2422.737  # M134: -2450.2104 -> -2422.737
$ws.Cells.Item(136, 9).Value = 6449.8335  # I136: 5742.7144 -> 6449.8335
$ws.Cells.Item(136, 10).Value = 1599.875  # J136: 1614.1428 -> 1599.875
$ws.Cells.Item(136, 11).Value = 19349.5005  # K136: 17228.1432 -> 19349.5005
$ws.Cells.Item(136, 12).Value = 4799.625  # L136: 4842.428400000001 -> 4799.625
$ws.Cells.Item(136, 13).Value = -16799.5005  # M136: -14678.1432 -> -16799.5005
$ws.Cells.Item(136, 14).Value = -9899.625  # N136: -9942.428400000001 -> -9899.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 556.7646999999999  # H5: 541.06665 -> 556.7646999999999
$ws.Cells.Item(5, 9).Value = 501.2  # I5: 508.5 -> 501.2
$ws.Cells.Item(5, 10).Value = 973.5  # J5: 997 -> 973.5
$ws.Cells.Item(5, 11).Value = 1503.6  # K5: 1525.5 -> 1503.6
$ws.Cells.Item(5, 12).Value = 2920.5  # L5: 2991 -> 2920.5
$ws.Cells.Item(5, 13).Value = -1391.6  # M5: -1413.5 -> -1391.6
$ws.Cells.Item(5, 14).Value = -3144.5  # N5: -3215 -> -3144.5
$ws.Cells.Item(131, 8).Value = 2560.3333  # H131: 2502 -> 2560.3333
$ws.Cells.Item(131, 9).Value = 1284.8  # I131: 1437.3334 -> 1284.8
$ws.Cells.Item(131, 10).Value = 3471.4285  # J131: 3566.6667 -> 3471.4285
$ws.Cells.Item(131, 11).Value = 3854.4  # K131: 4312.0002 -> 3854.4
$ws.Cells.Item(131, 12).Value = 10414.2855  # L131: 10700.0001 -> 10414.2855
$ws.Cells.Item(131, 13).Value = 1185.6  # M131: 727.9997999999996 -> 1185.6
$ws.Cells.Item(131, 14).Value = -20494.2855  # N131: -20780.0001 -> -20494.2855
$ws.Cells.Item(135, 8).Value = 556.7646999999999  # H135: 541.06665 -> 556.7646999999999
$ws.Cells.Item(135, 9).Value = 501.2  # I135: 508.5 -> 501.2
$ws.Cells.Item(135, 10).Value = 973.5  # J135: 997 -> 973.5
$ws.Cells.Item(135, 11).Value = 4510.8  # K135: 4576.5 -> 4510.8
$ws.Cells.Item(135, 12).Value = 8761.5  # L135: 8973 -> 8761.5
$ws.Cells.Item(135, 13).Value = -1975.8  # M135: -2041.5 -> -1975.8
$ws.Cells.Item(135, 14).Value = -13831.5  # N135: -14043 -> -13831.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6999  # H70: 5999.5 -> 6999
$ws.Cells.Item(70, 9).Value = 6999  # I70: 5999.5 -> 6999
$ws.Cells.Item(70, 11).Value = 6999  # K70: 5999.5 -> 6999
$ws.Cells.Item(70, 13).Value = -6729  # M70: -5729.5 -> -6729
$ws.Cells.Item(73, 8).Value = 6999  # H73: 5999.5 -> 6999
$ws.Cells.Item(73, 9).Value = 6999  # I73: 5999.5 -> 6999
$ws.Cells.Item(73, 11).Value = 6999  # K73: 5999.5 -> 6999
$ws.Cells.Item(73, 13).Value = -6063  # M73: -5063.5 -> -6063
$ws.Cells.Item(80, 8).Value = 3845.923  # H80: 4166.4165 -> 3845.923
$ws.Cells.Item(80, 9).Value = 2499.25  # I80: 2665.6667 -> 2499.25
$ws.Cells.Item(80, 10).Value = 4444.4443  # J80: 4666.6665 -> 4444.4443
$ws.Cells.Item(80, 11).Value = 2499.25  # K80: 2665.6667 -> 2499.25
$ws.Cells.Item(80, 12).Value = 4444.4443  # L80: 4666.6665 -> 4444.4443
$ws.Cells.Item(80, 13).Value = -1501.25  # M80: -1667.6667 -> -1501.25
$ws.Cells.Item(80, 14).Value = -6440.4443  # N80: -6662.6665 -> -6440.4443
$ws.Cells.Item(83, 8).Value = 3845.923  # H83: 4166.4165 -> 3845.923
$ws.Cells.Item(83, 9).Value = 2499.25  # I83: 2665.6667 -> 2499.25
$ws.Cells.Item(83, 10).Value = 4444.4443  # J83: 4666.6665 -> 4444.4443
$ws.Cells.Item(83, 11).Value = 12496.25  # K83: 13328.3335 -> 12496.25
$ws.Cells.Item(83, 12).Value = 22222.2215  # L83: 23333.3325 -> 22222.2215
$ws.Cells.Item(83, 13).Value = -7504.25  # M83: -8336.333500000001 -> -7504.25
$ws.Cells.Item(83, 14).Value = -32206.2215  # N83: -33317.3325 -> -32206.2215
$ws.Cells.Item(102, 8).Value = 1906.2307  # H102: 2042.4546 -> 1906.2307
$ws.Cells.Item(102, 9).Value = 1906.2307  # I102: 2042.4546 -> 1906.2307
$ws.Cells.Item(102, 11).Value = 1906.2307  # K102: 2042.4546 -> 1906.2307
$ws.Cells.Item(102, 13).Value = -284.2307000000001  # M102: -420.4546 -> -284.2307000000001
$ws.Cells.Item(122, 8).Value = 41512.445  # H122: 41534.39 -> 41512.445
$ws.Cells.Item(122, 9).Value = 40481.668  # I122: 40508 -> 40481.668
$ws.Cells.Item(122, 11).Value = 121445.004  # K122: 121524 -> 121445.004
$ws.Cells.Item(122, 13).Value = -118995.004  # M122: -119074 -> -118995.004
$ws.Cells.Item(126, 8).Value = 1216  # H126: 0 -> 1216
$ws.Cells.Item(126, 9).Value = 1250  # I126: 0 -> 1250
$ws.Cells.Item(126, 10).Value = 1114  # J126: 0 -> 1114
$ws.Cells.Item(126, 11).Value = 3750  # K126: 0 -> 3750
$ws.Cells.Item(126, 12).Value = 3342  # L126: 0 -> 3342
$ws.Cells.Item(126, 13).Value = -1280  # M126: (none) -> -1280
$ws.Cells.Item(126, 14).Value = -8282  # N126: (none) -> -8282
$ws.Cells.Item(132, 8).Value = 1333.3334  # H132: 1375 -> 1333.3334
$ws.Cells.Item(132, 10).Value = 2000  # J132: 2500 -> 2000
$ws.Cells.Item(132, 12).Value = 6000  # L132: 7500 -> 6000
$ws.Cells.Item(132, 14).Value = -11060  # N132: -12560 -> -11060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 778.25  # H9: 911.2222 -> 778.25
$ws.Cells.Item(9, 10).Value = 1350  # J9: 1506.25 -> 1350
$ws.Cells.Item(9, 12).Value = 1350  # L9: 1506.25 -> 1350
$ws.Cells.Item(9, 14).Value = -1798  # N9: -1954.25 -> -1798
$ws.Cells.Item(46, 8).Value = 1535.7368  # H46: 1552.1538 -> 1535.7368
$ws.Cells.Item(46, 9).Value = 862.63635  # I46: 899 -> 862.63635
$ws.Cells.Item(46, 10).Value = 2461.25  # J46: 3729.3333 -> 2461.25
$ws.Cells.Item(46, 11).Value = 862.63635  # K46: 899 -> 862.63635
$ws.Cells.Item(46, 12).Value = 2461.25  # L46: 3729.3333 -> 2461.25
$ws.Cells.Item(46, 13).Value = -674.63635  # M46: -711 -> -674.63635
$ws.Cells.Item(46, 14).Value = -2837.25  # N46: -4105.3333 -> -2837.25
$ws.Cells.Item(122, 8).Value = 4633.933  # H122: 4893.5 -> 4633.933
$ws.Cells.Item(122, 10).Value = 6250  # J122: 7300 -> 6250
$ws.Cells.Item(122, 12).Value = 18750  # L122: 21900 -> 18750
$ws.Cells.Item(122, 14).Value = -23650  # N122: -26800 -> -23650
$ws.Cells.Item(136, 8).Value = 3606.2222  # H136: 3580.4 -> 3606.2222
$ws.Cells.Item(136, 9).Value = 3493.875  # I136: 3477.6667 -> 3493.875
$ws.Cells.Item(136, 11).Value = 10481.625  # K136: 10433.0001 -> 10481.625
$ws.Cells.Item(136, 13).Value = -7931.625  # M136: -7883.000100000001 -> -7931.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 35000  # H49: 5000 -> 35000
$ws.Cells.Item(49, 9).Value = 0  # I49: 5000 -> 0
$ws.Cells.Item(49, 10).Value = 35000  # J49: 0 -> 35000
$ws.Cells.Item(49, 11).Value = 0  # K49: 5000 -> 0
$ws.Cells.Item(49, 12).Value = 35000  # L49: 0 -> 35000
$ws.Cells.Item(49, 13).ClearContents()  # M49: -4770 -> (cleared)
$ws.Cells.Item(49, 14).Value = -35460  # N49: (none) -> -35460
$ws.Cells.Item(101, 8).Value = 15300.25  # H101: 16867 -> 15300.25
$ws.Cells.Item(101, 10).Value = 15300.25  # J101: 16867 -> 15300.25
$ws.Cells.Item(101, 12).Value = 15300.25  # L101: 16867 -> 15300.25
$ws.Cells.Item(101, 14).Value = -21790.25  # N101: -23357 -> -21790.25
$ws.Cells.Item(122, 8).Value = 1563.5555  # H122: 1774.1428 -> 1563.5555
$ws.Cells.Item(122, 9).Value = 1563.5555  # I122: 1774.1428 -> 1563.5555
$ws.Cells.Item(122, 11).Value = 4690.666499999999  # K122: 5322.428400000001 -> 4690.666499999999
$ws.Cells.Item(122, 13).Value = -2240.666499999999  # M122: -2872.428400000001 -> -2240.666499999999
